$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (issue number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# --- Weekly crime statistics grid (rows 14-30) ---

# Row 14
$ws.Range("M14").Value = -44.444444444444

# Row 15
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 11
$ws.Range("K15").Value = -31.25
$ws.Range("M15").Value = -31.25
$ws.Range("N15").Value = -69.444444444444

# Row 16
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 162
$ws.Range("J16").Value = 174
$ws.Range("K16").Value = -6.896551724137
$ws.Range("L16").Value = 55.769230769230
$ws.Range("M16").Value = -37.451737451737
$ws.Range("N16").Value = -84.438040345821

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 12.5
$ws.Range("I17").Value = 355
$ws.Range("J17").Value = 364
$ws.Range("K17").Value = -2.472527472527
$ws.Range("L17").Value = 4.719764011799
$ws.Range("M17").Value = 71.497584541062
$ws.Range("N17").Value = -50.762829403606

# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 500
$ws.Range("F18").Value = 11
$ws.Range("H18").Value = -15.384615384615
$ws.Range("I18").Value = 106
$ws.Range("J18").Value = 157
$ws.Range("K18").Value = -32.484076433121
$ws.Range("L18").Value = -35.757575757575
$ws.Range("M18").Value = -32.484076433121
$ws.Range("N18").Value = -89.004149377593

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -30.769230769230
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 445
$ws.Range("J19").Value = 513
$ws.Range("K19").Value = -13.255360623781
$ws.Range("L19").Value = 8.009708737864
$ws.Range("M19").Value = -16.977611940298
$ws.Range("N19").Value = -33.582089552238

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 109
$ws.Range("J20").Value = 118
$ws.Range("K20").Value = -7.627118644067
$ws.Range("L20").Value = 55.714285714285
$ws.Range("M20").Value = -0.909090909090
$ws.Range("N20").Value = -90.404929577464

# Row 21
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 107
$ws.Range("H21").Value = -15.887850467289
$ws.Range("I21").Value = 1193
$ws.Range("J21").Value = 1347
$ws.Range("K21").Value = -11.432813659985
$ws.Range("L21").Value = 7.380738073807
$ws.Range("M21").Value = -7.805255023183
$ws.Range("N21").Value = -73.997384481255

# Row 22
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = -40
$ws.Range("M22").Value = -41.935483870967

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = -28.571428571428
$ws.Range("I23").Value = 160
$ws.Range("J23").Value = 173
$ws.Range("K23").Value = -7.514450867052
$ws.Range("L23").Value = -10.112359550561
$ws.Range("M23").Value = 35.593220338983

# Row 24
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 19.047619047619
$ws.Range("F24").Value = 82
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = -1.204819277108
$ws.Range("I24").Value = 969
$ws.Range("J24").Value = 1178
$ws.Range("K24").Value = -17.741935483871
$ws.Range("L24").Value = 14.134275618374
$ws.Range("M24").Value = -15.445026178010

# Row 25
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -9.090909090909
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 20.588235294117
$ws.Range("I25").Value = 589
$ws.Range("J25").Value = 540
$ws.Range("K25").Value = 9.074074074074
$ws.Range("L25").Value = 31.180400890868
$ws.Range("M25").Value = 15.717092337917

# Row 26
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 19
$ws.Range("K26").Value = -32.142857142857
$ws.Range("L26").Value = -17.391304347826

# Row 27
$ws.Range("I27").Value = 51
$ws.Range("K27").Value = -5.555555555555
$ws.Range("L27").Value = 34.210526315789

# Row 28
$ws.Range("F28").Value = 4
$ws.Range("I28").Value = 12
$ws.Range("K28").Value = -53.846153846153
$ws.Range("L28").Value = -33.333333333333
$ws.Range("M28").Value = -47.826086956521
$ws.Range("N28").Value = -82.352941176470

# Row 29
$ws.Range("F29").Value = 3
$ws.Range("I29").Value = 8
$ws.Range("K29").Value = -52.941176470588
$ws.Range("L29").Value = -46.666666666666
$ws.Range("M29").Value = -55.555555555555
$ws.Range("N29").Value = -84.905660377358

# Row 30
$ws.Range("I30").Value = 12
$ws.Range("K30").Value = 71.428571428571
$ws.Range("L30").Value = 100

# --- Cells that changed from a text placeholder ("0" / "***.*") to a real number ---
# need an explicit number format so the cell is stored as numeric (integer / one-decimal) style,
# matching the rest of the data grid, instead of inheriting the old text format.
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
